$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Change the scenario dropdown selection in C11
$ws.Range("C11").Value = "2 Slightly over limit at min ded level"

# 2. Update the max-ded formula in C32 (and the shared formula across D32:E32)
$ws.Range("C32").Formula = "=MIN(C29,C27-C30,C17-C30)"
$ws.Range("D32:E32").Formula = "=MIN(D29,D27-D30,D17-D30)"

# 3. Apply number formatting to D25 (center aligned, integer thousands format)
$ws.Range("D25").NumberFormat = "#,##0"

# 4. Update the active selection to C36:D36
$ws.Range("C36:D36").Select()
